$wb = $excel.ActiveWorkbook

# --- general sheet ---
$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Activate()

# B4: was formula =B5*0.5 -> becomes a plain value 550
$wsGeneral.Range("B4").Value = 550

# B46: 1.6242000000000001 -> 0
$wsGeneral.Range("B46").Value = 0

# B51: -1738 -> 1738 ; E51 formula changes from shared "=B51" to "=B51*1000"
$wsGeneral.Range("B51").Value = 1738
$wsGeneral.Range("E51").Formula = "=B51*1000"

# Selection / active view for general
$wsGeneral.Range("E52").Select()

# --- initialConditions sheet ---
$wsInit = $wb.Worksheets.Item("initialConditions")
$wsInit.Activate()

# B3: 10.3957 (value) -> formula =10.3957+1737.5
$wsInit.Range("B3").Formula = "=10.3957+1737.5"
$wsInit.Range("B3").Borders.LineStyle = -4142

# B7: -0.58329957806655397 -> 0
$wsInit.Range("B7").Value = 0

# B8: 1 -> 1E-3
$wsInit.Range("B8").Value = 0.001

# B12: 1 -> 1E-3 ; C12 gets unit label "m" ; E12 formula becomes standalone "=B12*1000"
$wsInit.Range("B12").Value = 0.001
$wsInit.Range("C12").Value = "m"
$wsInit.Range("E12").Formula = "=B12*1000"

$wsInit.Range("B13").Select()

# --- navStateIdx sheet: no longer the active tab ---
$wsNavIdx = $wb.Worksheets.Item("navStateIdx")
$wsNavIdx.Range("A7").Select()

# --- truthStateParams sheet: selection change only ---
$wsTruthParams = $wb.Worksheets.Item("truthStateParams")
$wsTruthParams.Range("A15:C19").Select()

# --- errorInjection sheet ---
$wsErr = $wb.Worksheets.Item("errorInjection")
$wsErr.Activate()

$wsErr.Range("B8").Value = 10
$wsErr.Range("B9").Value = 1
$wsErr.Range("B10").Value = 2
$wsErr.Range("B11").Value = 3
$wsErr.Range("B12").Value = 10
$wsErr.Range("B13").Value = 1
$wsErr.Range("B14").Value = 1
$wsErr.Range("B15").Value = 1

$wsErr.Range("B16").Select()

# Re-activate "general" last so it becomes the active tab (tabSelected) in the saved file
$wsGeneral.Activate()
$wsGeneral.Range("E52").Select()
